# Update the position/size of shape "object 11" (the "Gerente de sistemas"
# label box) on slide 1, and mark its trailing paragraph-mark runs as
# "clean" (dirty="0") same as a real PowerPoint resize/reflow would.
#
# Target OOXML (from the diff):
#   <a:off x="6067424" y="2189479"/>
#   <a:ext cx="638175" cy="269304"/>
#   ... <a:endParaRPr sz="900" dirty="0"> (x2)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape "object 11" (cNvPr id="11") is the 10th item of the slide's Shapes
# collection (the very first spTree child - the group's own nvGrpSpPr - is
# not exposed as Shapes.Item(1), so ids are offset by one from the 1-based
# COM index).
$shp = $s.Shapes.Item(10)

# Sanity checks so we fail loudly instead of silently editing the wrong shape.
if ($shp.Name -ne "object 11") {
    throw "Expected shape 'object 11' at index 10, found '$($shp.Name)'"
}

# Reposition / resize the textbox. PowerPoint's COM surface works in points
# (1 pt = 12700 EMU); the literals below are chosen so the round-tripped EMU
# values exactly match the target off/ext from the diff:
#   x  -> 6067424 EMU
#   y  -> 2189479 EMU (left untouched - unchanged by the edit)
#   cx -> 638175  EMU
#   cy -> 269304  EMU
$shp.Left = 477.7499562598425
$shp.Width = 50.25
$shp.Height = 21.20503937007874

# The text itself ("Gerente de" / "sistemas") is unchanged; only the
# paragraph-end run properties picked up a dirty="0" flag (PowerPoint marks
# them "not dirty for spell-check" after the edit/reflow above).
$tf = $shp.TextFrame
$tr = $tf.TextRange
for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    $para = $tr.Paragraphs($i)
    $para.Font.Size = $para.Font.Size
}
